$d = $word.ActiveDocument

# --- Change 1: fix wording of the "mucosa" sentence -----------------------
$d.Content.Find.Execute(
    "Cancers start on the very inside of the layer called the mucosa",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cancers start on the very inside layer called the mucosa", 2) | Out-Null

# --- Change 2: split the "Chemotherapy + radiation ..." paragraph ---------
# A new "FirstParagraph"-styled paragraph carrying the patient-count
# sentence is inserted immediately before it; the original paragraph
# (content/runs unchanged) is then demoted to "BodyText" style.
$rng = $d.Content
$rng.Find.Execute(
    "Chemotherapy + radiation given together over 6 weeks",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$targetPara = $rng.Paragraphs(1)
$insertionPoint = $targetPara.Range.Start

# Inserting before reseats any Paragraph object anchored at
# $insertionPoint onto the freshly-created (empty) paragraph, so work off
# Range positions rather than the old $targetPara reference from here on.
$targetPara.Range.InsertParagraphBefore() | Out-Null

$newPara = $d.Range($insertionPoint, $insertionPoint).Paragraphs(1)
$newPara.Range.Text = "363 patients with esophageal cancer studied"

$afterNewParaEnd = $newPara.Range.End
$shiftedTargetPara = $d.Range($afterNewParaEnd, $afterNewParaEnd).Paragraphs(1)
$shiftedTargetPara.Style = "BodyText"
